$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# "These teams ... between their time ... college basketball" used to be
# split across three runs (with spellcheck proofErr markers bracketing
# "their "). The revision removes the proofErr markers and folds the
# text back into a single run. We reproduce that by replacing the whole
# span with itself; the host's Find/Replace machinery re-serialises the
# touched runs (dropping the now-unneeded proofErr markers).
$oldText1 = "These teams are often the best so they must be the best programs at developing players between their time from high school prospects to professionals NBA players, one might believe. A quick look at the top of the draft boards and we will find it littered with players from the top programs in college basketball"
$found1 = $d.Content.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, $oldText1, 2)

# --- Edit 2 -----------------------------------------------------------
# "But do these teams ... trajectory to the NBA.  " (two trailing
# spaces) becomes "...trajectory to the NBA. " (one trailing space)
# followed immediately by a new sentence about the analysis' focus.
$apos = [char]0x2019
$oldText2 = "But do these teams really develop the best players or do they recruit the best players who were already on the trajectory to the NBA.  "
$newText2 = "But do these teams really develop the best players or do they recruit the best players who were already on the trajectory to the NBA. For this analysis, we will focus on the top 100 ranked high school players from each year, and we" + $apos + "ll see how many of them get drafted into the NBA and whether the college that they attend will affect their chances. "
$found2 = $d.Content.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, $newText2, 2)

# --- New paragraphs -----------------------------------------------------
# A blank paragraph, then a short intro line, a blank, a "GRAPH"
# placeholder, a blank, a paragraph of commentary, a blank, another
# commentary paragraph, a blank, another "GRAPH" placeholder, a blank,
# another commentary paragraph, a blank, the closing commentary
# paragraph, and finally four blank paragraphs.

# empty paragraph
$null = $d.Paragraphs.Add()

$p = $d.Paragraphs.Add()
$p.Range.InsertAfter("Let" + $apos + "s look at how often top 100 prospects get drafted into the NBA based on the school they go to.  ")

$null = $d.Paragraphs.Add()

$p = $d.Paragraphs.Add()
$p.Range.InsertAfter("GRAPH")

$null = $d.Paragraphs.Add()

$p = $d.Paragraphs.Add()
$p.Range.InsertAfter("So, we can see that the colleges that recruit better high school prospects, have more of them drafted into the NBA. The schools on the left are the top programs in the country such as Duke Kentucky and UNC. These schools recruit highly ranked prospects and have a high draft rate. ")

$null = $d.Paragraphs.Add()

$p = $d.Paragraphs.Add()
$p.Range.InsertAfter("A prospect" + $apos + "s high school ranking will affect their chances of being drafted. Can a school be credited for recruiting a good prospect and having that prospect end up in the NBA or was that prospect going to make it no matter what school he went to? To evaluate whether a school increased a player" + $apos + "s chance of being drafted we need to compare how often a school" + $apos + "s players get drafted to how often a player of the same ranking is usually drafted. ")

$null = $d.Paragraphs.Add()

$p = $d.Paragraphs.Add()
$p.Range.InsertAfter("GRAPH")

$null = $d.Paragraphs.Add()

$p = $d.Paragraphs.Add()
$p.Range.InsertAfter("The odds of a player being drafted to the NBA decrease as their ranking gets lower. Every prospect ranked 1 or 2 have been drafted into the NBA. There is some noise in the data so the actual chances to being drafted are modeled with the overlayed line. Interesting to see how tough it is to make it to the NBA. Anyone ranked outside the top 20 in high school rankings has an outside chance of being drafted. ")

$null = $d.Paragraphs.Add()

$p = $d.Paragraphs.Add()
$p.Range.InsertAfter("Using these theoretical odds of being drafted, we can calculate how much each school increased their player" + $apos + "s chances of being drafted. This is turned into a draft score for each school. The higher the average draft score the better.")

$null = $d.Paragraphs.Add()
$null = $d.Paragraphs.Add()
$null = $d.Paragraphs.Add()
$null = $d.Paragraphs.Add()

Write-Host "found1=$found1 found2=$found2 paraCount=$($d.Paragraphs.Count)"
